# Apply "adding new progress as of date 04 nov 2025" update
# On the "Training Dashboard" worksheet, for rows 3 through 16:
#   - decrement column H (PERIOD TO EXPIRE) by 1
#   - update column I (LAST UPDATE) text from "03-Nov-2025" to "04-Nov-2025"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 16; $row++) {
    # column H (PERIOD TO EXPIRE) - decrease by 1 day
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value = $hCell.Value2 - 1

    # column I (LAST UPDATE) - move the progress date forward to 04-Nov-2025
    # The leading apostrophe keeps Excel from reinterpreting the literal
    # text as a real date value (the column is stored as plain text).
    $iCell = $ws.Cells.Item($row, 9)
    $iCell.Value2 = "'04-Nov-2025"
}
